$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column from 2023-11-13 (45243) to 2023-11-14 (45244)
# for rows 2 through 5, preserving the existing date formatting on the cells.
foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
